# Adds five new per-contract metric columns (H:L) to the "Contracts" sheet:
#   H = Number of structs
#   I = Number of using-for
#   J = Number of custom error definitions
#   K = Number of events
#   L = Number of inherited classes
# This mirrors the upstream diff that extends the sheet dimension from A1:G77 to A1:L77.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1), columns H-L ---
$ws.Range("H1").Value = "Number of structs"
$ws.Range("I1").Value = "Number of using-for"
$ws.Range("J1").Value = "Number of custom error definitions"
$ws.Range("K1").Value = "Number of events"
$ws.Range("L1").Value = "Number of inherited classes"

# --- Bulk-fill the new metric values for data rows 2-77 (H2:L77) via one 2D array write ---
$data = New-Object "object[,]" 76,5
$data[0,0] = 2
$data[0,1] = 6
$data[0,2] = 4
$data[0,3] = 2
$data[0,4] = 2
$data[1,0] = 4
$data[1,1] = 6
$data[1,2] = 7
$data[1,3] = 4
$data[1,4] = 2
$data[2,0] = 0
$data[2,1] = 2
$data[2,2] = 3
$data[2,3] = 6
$data[2,4] = 2
$data[3,0] = 1
$data[3,1] = 2
$data[3,2] = 0
$data[3,3] = 1
$data[3,4] = 2
$data[4,0] = 0
$data[4,1] = 0
$data[4,2] = 0
$data[4,3] = 4
$data[4,4] = 0
$data[5,0] = 1
$data[5,1] = 0
$data[5,2] = 0
$data[5,3] = 7
$data[5,4] = 1
$data[6,0] = 0
$data[6,1] = 0
$data[6,2] = 0
$data[6,3] = 3
$data[6,4] = 2
$data[7,0] = 0
$data[7,1] = 0
$data[7,2] = 0
$data[7,3] = 2
$data[7,4] = 0
$data[8,0] = 0
$data[8,1] = 0
$data[8,2] = 0
$data[8,3] = 13
$data[8,4] = 0
$data[9,0] = 0
$data[9,1] = 0
$data[9,2] = 0
$data[9,3] = 10
$data[9,4] = 0
$data[10,0] = 0
$data[10,1] = 0
$data[10,2] = 0
$data[10,3] = 0
$data[10,4] = 0
$data[11,0] = 0
$data[11,1] = 0
$data[11,2] = 0
$data[11,3] = 0
$data[11,4] = 0
$data[12,0] = 0
$data[12,1] = 0
$data[12,2] = 0
$data[12,3] = 2
$data[12,4] = 1
$data[13,0] = 0
$data[13,1] = 0
$data[13,2] = 0
$data[13,3] = 0
$data[13,4] = 0
$data[14,0] = 0
$data[14,1] = 0
$data[14,2] = 0
$data[14,3] = 0
$data[14,4] = 0
$data[15,0] = 0
$data[15,1] = 0
$data[15,2] = 0
$data[15,3] = 0
$data[15,4] = 0
$data[16,0] = 0
$data[16,1] = 0
$data[16,2] = 0
$data[16,3] = 0
$data[16,4] = 0
$data[17,0] = 0
$data[17,1] = 0
$data[17,2] = 0
$data[17,3] = 0
$data[17,4] = 0
$data[18,0] = 0
$data[18,1] = 0
$data[18,2] = 0
$data[18,3] = 0
$data[18,4] = 0
$data[19,0] = 0
$data[19,1] = 0
$data[19,2] = 0
$data[19,3] = 0
$data[19,4] = 0
$data[20,0] = 0
$data[20,1] = 0
$data[20,2] = 0
$data[20,3] = 0
$data[20,4] = 0
$data[21,0] = 0
$data[21,1] = 0
$data[21,2] = 0
$data[21,3] = 0
$data[21,4] = 0
$data[22,0] = 0
$data[22,1] = 1
$data[22,2] = 0
$data[22,3] = 0
$data[22,4] = 1
$data[23,0] = 0
$data[23,1] = 0
$data[23,2] = 0
$data[23,3] = 0
$data[23,4] = 0
$data[24,0] = 0
$data[24,1] = 2
$data[24,2] = 0
$data[24,3] = 0
$data[24,4] = 1
$data[25,0] = 0
$data[25,1] = 0
$data[25,2] = 0
$data[25,3] = 0
$data[25,4] = 1
$data[26,0] = 0
$data[26,1] = 2
$data[26,2] = 0
$data[26,3] = 0
$data[26,4] = 1
$data[27,0] = 0
$data[27,1] = 2
$data[27,2] = 0
$data[27,3] = 0
$data[27,4] = 1
$data[28,0] = 0
$data[28,1] = 4
$data[28,2] = 0
$data[28,3] = 0
$data[28,4] = 1
$data[29,0] = 3
$data[29,1] = 0
$data[29,2] = 0
$data[29,3] = 0
$data[29,4] = 0
$data[30,0] = 0
$data[30,1] = 0
$data[30,2] = 0
$data[30,3] = 0
$data[30,4] = 0
$data[31,0] = 0
$data[31,1] = 0
$data[31,2] = 0
$data[31,3] = 0
$data[31,4] = 0
$data[32,0] = 0
$data[32,1] = 0
$data[32,2] = 0
$data[32,3] = 0
$data[32,4] = 0
$data[33,0] = 3
$data[33,1] = 2
$data[33,2] = 0
$data[33,3] = 0
$data[33,4] = 0
$data[34,0] = 11
$data[34,1] = 0
$data[34,2] = 0
$data[34,3] = 0
$data[34,4] = 0
$data[35,0] = 2
$data[35,1] = 2
$data[35,2] = 0
$data[35,3] = 2
$data[35,4] = 1
$data[36,0] = 0
$data[36,1] = 3
$data[36,2] = 1
$data[36,3] = 1
$data[36,4] = 1
$data[37,0] = 0
$data[37,1] = 5
$data[37,2] = 6
$data[37,3] = 22
$data[37,4] = 1
$data[38,0] = 0
$data[38,1] = 0
$data[38,2] = 0
$data[38,3] = 0
$data[38,4] = 2
$data[39,0] = 0
$data[39,1] = 7
$data[39,2] = 1
$data[39,3] = 0
$data[39,4] = 1
$data[40,0] = 0
$data[40,1] = 3
$data[40,2] = 2
$data[40,3] = 3
$data[40,4] = 1
$data[41,0] = 0
$data[41,1] = 1
$data[41,2] = 2
$data[41,3] = 3
$data[41,4] = 1
$data[42,0] = 0
$data[42,1] = 0
$data[42,2] = 0
$data[42,3] = 0
$data[42,4] = 1
$data[43,0] = 0
$data[43,1] = 2
$data[43,2] = 3
$data[43,3] = 6
$data[43,4] = 2
$data[44,0] = 1
$data[44,1] = 2
$data[44,2] = 0
$data[44,3] = 1
$data[44,4] = 2
$data[45,0] = 0
$data[45,1] = 0
$data[45,2] = 0
$data[45,3] = 0
$data[45,4] = 0
$data[46,0] = 1
$data[46,1] = 0
$data[46,2] = 0
$data[46,3] = 0
$data[46,4] = 0
$data[47,0] = 1
$data[47,1] = 0
$data[47,2] = 0
$data[47,3] = 0
$data[47,4] = 0
$data[48,0] = 1
$data[48,1] = 0
$data[48,2] = 0
$data[48,3] = 0
$data[48,4] = 0
$data[49,0] = 1
$data[49,1] = 0
$data[49,2] = 0
$data[49,3] = 0
$data[49,4] = 1
$data[50,0] = 1
$data[50,1] = 0
$data[50,2] = 0
$data[50,3] = 0
$data[50,4] = 1
$data[51,0] = 0
$data[51,1] = 0
$data[51,2] = 0
$data[51,3] = 0
$data[51,4] = 0
$data[52,0] = 0
$data[52,1] = 0
$data[52,2] = 0
$data[52,3] = 0
$data[52,4] = 0
$data[53,0] = 0
$data[53,1] = 0
$data[53,2] = 0
$data[53,3] = 0
$data[53,4] = 0
$data[54,0] = 0
$data[54,1] = 0
$data[54,2] = 0
$data[54,3] = 0
$data[54,4] = 0
$data[55,0] = 0
$data[55,1] = 0
$data[55,2] = 0
$data[55,3] = 0
$data[55,4] = 0
$data[56,0] = 0
$data[56,1] = 0
$data[56,2] = 0
$data[56,3] = 0
$data[56,4] = 0
$data[57,0] = 0
$data[57,1] = 0
$data[57,2] = 0
$data[57,3] = 0
$data[57,4] = 0
$data[58,0] = 0
$data[58,1] = 1
$data[58,2] = 0
$data[58,3] = 0
$data[58,4] = 1
$data[59,0] = 0
$data[59,1] = 0
$data[59,2] = 0
$data[59,3] = 0
$data[59,4] = 0
$data[60,0] = 0
$data[60,1] = 1
$data[60,2] = 0
$data[60,3] = 0
$data[60,4] = 1
$data[61,0] = 0
$data[61,1] = 0
$data[61,2] = 0
$data[61,3] = 0
$data[61,4] = 2
$data[62,0] = 0
$data[62,1] = 1
$data[62,2] = 0
$data[62,3] = 0
$data[62,4] = 1
$data[63,0] = 0
$data[63,1] = 2
$data[63,2] = 0
$data[63,3] = 0
$data[63,4] = 1
$data[64,0] = 0
$data[64,1] = 1
$data[64,2] = 1
$data[64,3] = 0
$data[64,4] = 1
$data[65,0] = 0
$data[65,1] = 2
$data[65,2] = 1
$data[65,3] = 0
$data[65,4] = 1
$data[66,0] = 0
$data[66,1] = 0
$data[66,2] = 2
$data[66,3] = 0
$data[66,4] = 0
$data[67,0] = 4
$data[67,1] = 2
$data[67,2] = 0
$data[67,3] = 0
$data[67,4] = 0
$data[68,0] = 11
$data[68,1] = 0
$data[68,2] = 0
$data[68,3] = 0
$data[68,4] = 0
$data[69,0] = 2
$data[69,1] = 2
$data[69,2] = 0
$data[69,3] = 2
$data[69,4] = 1
$data[70,0] = 0
$data[70,1] = 2
$data[70,2] = 1
$data[70,3] = 1
$data[70,4] = 1
$data[71,0] = 0
$data[71,1] = 2
$data[71,2] = 5
$data[71,3] = 21
$data[71,4] = 1
$data[72,0] = 0
$data[72,1] = 0
$data[72,2] = 0
$data[72,3] = 0
$data[72,4] = 2
$data[73,0] = 0
$data[73,1] = 3
$data[73,2] = 2
$data[73,3] = 0
$data[73,4] = 1
$data[74,0] = 4
$data[74,1] = 3
$data[74,2] = 19
$data[74,3] = 9
$data[74,4] = 2
$data[75,0] = 0
$data[75,1] = 1
$data[75,2] = 2
$data[75,3] = 0
$data[75,4] = 2

$ws.Range("H2:L77").Value = $data

# --- Best-effort: extend the "numberStoredAsText" ignored-error range to the new A1:L77 extent.
#     (Cosmetic worksheet metadata; harmless no-op if this COM surface does not expose it.)
try {
    $ws.Range("A1:L77").Errors.Item(9).Ignore = $true
} catch {
}
try {
    $ws.Range("A1:L77").ErrorCheckingOptions.NumberAsText = $true
} catch {
}

Write-Host ("Updated range now: " + $ws.UsedRange.Address())
